# Fruta / hortaliza, semanal
# A new weekly record was added to the "Coco" sheet. It becomes the new
# row 26 (most recent week), pushing all the former rows 26-59 down by
# one position (to 27-60). All other data stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 26 - this shifts the existing
# rows 26-59 down to 27-60, matching the rest of the diff.
$ws.Rows.Item(26).Insert()

# Populate the freshly inserted row 26 with the new weekly record.
$ws.Cells.Item(26, 1).Value  = 10
$ws.Cells.Item(26, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(26, 3).Value  = "La Araucanía"
$ws.Cells.Item(26, 4).Value  = 44671
$ws.Cells.Item(26, 5).Value  = 9
$ws.Cells.Item(26, 6).Value  = "Fruta"
$ws.Cells.Item(26, 7).Value  = 100108
$ws.Cells.Item(26, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(26, 9).Value  = 100108007
$ws.Cells.Item(26, 10).Value = "Coco"
$ws.Cells.Item(26, 11).Value = "Sin especificar"
$ws.Cells.Item(26, 12).Value = "Primera"
$ws.Cells.Item(26, 13).Value = 25
$ws.Cells.Item(26, 14).Value = 30000
$ws.Cells.Item(26, 15).Value = 30000
$ws.Cells.Item(26, 16).Value = 30000
$ws.Cells.Item(26, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(26, 18).Value = "Perú"
$ws.Cells.Item(26, 19).Value = 1500
$ws.Cells.Item(26, 20).Value = 20
